$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph with step "2." : remove redundant sz/szCs overrides ---
$p1 = Get-ParagraphByText $d "*funcionário insere no sistema o nome completo*"
if ($p1 -ne $null) {
    $body1 = '<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000037"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2.</w:t></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> O funcionário insere no sistema o nome completo, data de nascimento, número e tipo do documento de identificação do passageiro; e seleciona o voo desejado pelo cliente e a classe da passagem.</w:t></w:r></w:p>'
    $xml1 = $xmlHeader + $body1 + $xmlFooter
    $p1.Range.InsertXML($xml1)
} else {
    Write-Host "WARN: paragraph 2. not found"
}

# --- Paragraph with step "5." : remove redundant sz/szCs overrides + extend text ---
$p2 = Get-ParagraphByText $d "*O sistema exibe o preço*"
if ($p2 -ne $null) {
    $body2 = '<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000003A"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">5.</w:t></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> O sistema exibe o preço da passagem e oferece as opções de 25%, 50% ou 100% de desconto, de acordo com a quantidade de pontos possuídos pelo passageiro; além da opção sem desconto (0%), que não usa pontos.</w:t></w:r></w:p>'
    $xml2 = $xmlHeader + $body2 + $xmlFooter
    $p2.Range.InsertXML($xml2)
} else {
    Write-Host "WARN: paragraph 5. not found"
}

Write-Host "done"
